$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 376; existing rows 376-389 shift down to 377-390.
$ws.Rows.Item(376).Insert()

# Populate the newly inserted row 376 with the new weekly record.
$ws.Range("A376").Value = 10
$ws.Range("B376").Value = "Vega Modelo de Temuco"
$ws.Range("C376").Value = "La Araucanía"
$ws.Range("D376").Value = 45267
$ws.Range("D376").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E376").Value = 9
$ws.Range("F376").Value = 100114007
$ws.Range("G376").Value = "Jengibre"
$ws.Range("H376").Value = "Sin especificar"
$ws.Range("I376").Value = "Primera"
$ws.Range("J376").Value = 190
$ws.Range("K376").Value = 25000
$ws.Range("L376").Value = 26000
$ws.Range("M376").Value = 25474
$ws.Range("N376").Value = "$/caja 13 kilos"
$ws.Range("O376").Value = "Perú"
$ws.Range("P376").Value = 1960
$ws.Range("Q376").Value = 13
$ws.Range("R376").Value = "Hortaliza"
